# L5CG3_routine_3rdSem — updated for new routine
# Rebuild row 1 as a single title cell, and re-map the data columns
# (Module Code / Title / Hours shift right by one, Room/Block/Group
# re-ordered, and the old Level/Course columns are dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: collapse the header row into a single title cell ---
$ws.Range("B1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Data rows 2-10: new column layout ---
# Col A = Day (unchanged)
# Col B = Time (unchanged)
# Col C = Hours (was col E, numeric)
# Col D = Module Code (was col C)
# Col E = Module Title (was col D)
# Col F = Class Type (unchanged)
# Col G = Lecturer (unchanged)
# Col H = Group (was col J)
# Col I = Block (was col I, except row 3 which changes WLV -> HCK)
# Col J = Room (was col H)
# Col K, L (Level, Course) removed

$rows = @(
    @{ Row = 2;  Hours = 2.5; Code = "5CS020"; Title = "Distributed and Cloud Systems Programming"; Group = "L5CG3";          Block = "WLV"; Room = "Lab-02 Moseley" },
    @{ Row = 3;  Hours = 2;   Code = "5CS024"; Title = "Collaborative Development";                 Group = "L5CG3";          Block = "HCK"; Room = "TR-05 Ranipokhari" },
    @{ Row = 4;  Hours = 2.5; Code = "5CS024"; Title = "Collaborative Development";                 Group = "L5CG3";          Block = "WLV"; Room = "TR-01 Dudley" },
    @{ Row = 5;  Hours = 2;   Code = "5CS022"; Title = "Human Computer Interaction";                Group = "L5CG(1+2+3+4)";  Block = "WLV"; Room = "LT-02 Telford" },
    @{ Row = 6;  Hours = 2;   Code = "5CS020"; Title = "Distributed and Cloud Systems Programming"; Group = "L5CG(1+2+3+4)";  Block = "WLV"; Room = "LT-01 Wulfruna" },
    @{ Row = 7;  Hours = 2;   Code = "5CS024"; Title = "Collaborative Development";                 Group = "L5CG(1+2+3+4)";  Block = "WLV"; Room = "LT-01 Wulfruna" },
    @{ Row = 8;  Hours = 2;   Code = "5CS022"; Title = "Human Computer Interaction";                Group = "L5CG3";          Block = "WLV"; Room = "TR-02 Stafford" },
    @{ Row = 9;  Hours = 2;   Code = "5CS020"; Title = "Distributed and Cloud Systems Programming"; Group = "L5CG3";          Block = "WLV"; Room = "TR-03 Westbromwich" },
    @{ Row = 10; Hours = 2.5; Code = "5CS022"; Title = "Human Computer Interaction";                Group = "L5CG3";          Block = "WLV"; Room = "TR-03 Westbromwich" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 3).Value = $r.Hours    # C: Hours (numeric)
    $ws.Cells.Item($row, 4).Value = $r.Code     # D: Module Code
    $ws.Cells.Item($row, 5).Value = $r.Title    # E: Module Title
    $ws.Cells.Item($row, 8).Value = $r.Group    # H: Group
    $ws.Cells.Item($row, 9).Value = $r.Block    # I: Block
    $ws.Cells.Item($row, 10).Value = $r.Room    # J: Room
}

# Drop the old Level (K) and Course (L) columns entirely.
$ws.Range("K2:L10").ClearContents()
